# Updated symbol list on Tue Jan 24 04:38:46 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '321.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '5.54%'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '36.24'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.17%'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.118'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.72%'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08101'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.58%'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.152'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-1.14%'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '8.029'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.52%'

$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9291'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.19%'

$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1002'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '3.32%'

$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1883'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.31%'

$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09226'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '6.33%'

$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03588'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '3.17%'

$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09949'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.38%'

$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001439'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.43%'

$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005689'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.24%'

$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.465'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.14%'

$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.134'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.95%'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '18.01%'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3372'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-1.52%'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.48%'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.050'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '5.71%'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.20%'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04604'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.08%'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.07%'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004745'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-6.98%'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001300'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-7.13%'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0004502'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-5.21%'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02042'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '11.86%'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04997'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '5.72%'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007800'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2.20%'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1400'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.45%'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007810'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '0.80%'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002081'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-6.68%'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01214'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '9.31%'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006407'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.48%'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.02%'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '17.39%'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-4.95%'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002100'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.02%'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002000'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.02%'
